$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H86").Value = 7071.2856
$ws.Range("J86").Value = 8012.25
$ws.Range("L86").Value = 8012.25
$ws.Range("N86").Value = -10258.25
$ws.Range("H87").Value = 26665
$ws.Range("J87").Value = 26665
$ws.Range("L87").Value = 26665
$ws.Range("N87").Value = -29161
$ws.Range("H89").Value = 7071.2856
$ws.Range("J89").Value = 8012.25
$ws.Range("L89").Value = 40061.25
$ws.Range("N89").Value = -51293.25
$ws.Range("H90").Value = 26665
$ws.Range("J90").Value = 26665
$ws.Range("L90").Value = 79995
$ws.Range("N90").Value = -92475
$ws.Range("H111").Value = 949.5
$ws.Range("I111").Value = 900
$ws.Range("J111").Value = 999
$ws.Range("K111").Value = 2700
$ws.Range("L111").Value = 2997
$ws.Range("M111").Value = 367
$ws.Range("N111").Value = -9131
$ws.Range("H116").Value = 4670.3335
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 5005.5
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 5005.5
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -11889.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2257
$ws.Range("I2").Value = 705.8
$ws.Range("K2").Value = 705.8
$ws.Range("M2").Value = -592.8
$ws.Range("H116").Value = 2257
$ws.Range("I116").Value = 705.8
$ws.Range("K116").Value = 705.8
$ws.Range("M116").Value = 1588.2
$ws.Range("H132").Value = 2405
$ws.Range("I132").Value = 2405
$ws.Range("K132").Value = 7215
$ws.Range("M132").Value = -4685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2257
$ws.Range("I3").Value = 705.8
$ws.Range("K3").Value = 705.8
$ws.Range("M3").Value = -591.8
$ws.Range("H82").Value = 19917.445
$ws.Range("H85").Value = 19917.445
$ws.Range("H107").Value = 4998.5
$ws.Range("J107").Value = 4998.5
$ws.Range("L107").Value = 4998.5
$ws.Range("M107").Value = -8838.5
$ws.Range("H110").Value = 107188
$ws.Range("J110").Value = 107188
$ws.Range("L110").Value = 107188
$ws.Range("N110").Value = -115368
$ws.Range("H134").Value = 5776.2173
$ws.Range("I134").Value = 5197.864
$ws.Range("J134").Value = 18500
$ws.Range("K134").Value = 15593.592
$ws.Range("L134").Value = 55500
$ws.Range("M134").Value = -13058.592
$ws.Range("N134").Value = -60570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1397.125
$ws.Range("I31").Value = 1396.4
$ws.Range("K31").Value = 1396.4
$ws.Range("M31").Value = -1101.4
$ws.Range("H34").Value = 1397.125
$ws.Range("I34").Value = 1396.4
$ws.Range("K34").Value = 1396.4
$ws.Range("M34").Value = -1194.4
$ws.Range("H58").Value = 1518.8846
$ws.Range("I58").Value = 1602.6666
$ws.Range("J58").Value = 1447.0714
$ws.Range("K58").Value = 1602.6666
$ws.Range("L58").Value = 1447.0714
$ws.Range("M58").Value = -1399.6666
$ws.Range("N58").Value = -1853.0714
$ws.Range("H62").Value = 4166.6665
$ws.Range("I62").Value = 3750
$ws.Range("K62").Value = 3750
$ws.Range("M62").Value = -3126
$ws.Range("H65").Value = 4166.6665
$ws.Range("I65").Value = 3750
$ws.Range("K65").Value = 18750
$ws.Range("M65").Value = -15630
$ws.Range("H136").Value = 1518.8846
$ws.Range("I136").Value = 1602.6666
$ws.Range("J136").Value = 1447.0714
$ws.Range("K136").Value = 4807.9998
$ws.Range("L136").Value = 4341.2142
$ws.Range("M136").Value = -2257.9998
$ws.Range("N136").Value = -9441.2142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 250
$ws.Range("I49").Value = 250
$ws.Range("K49").Value = 750
$ws.Range("M49").Value = -594
$ws.Range("H56").Value = 10930.77
$ws.Range("I56").Value = 10930.77
$ws.Range("K56").Value = 10930.77
$ws.Range("M56").Value = -10400.77
$ws.Range("H96").Value = 2500
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 7500
$ws.Range("M96").Value = -5441
$ws.Range("H103").Value = 85.75
$ws.Range("J103").Value = 50
$ws.Range("L103").Value = 150
$ws.Range("N103").Value = -1908
$ws.Range("H132").Value = 1309.125
$ws.Range("I132").Value = 1371.75
$ws.Range("K132").Value = 12345.75
$ws.Range("M132").Value = -9815.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2521.5625
$ws.Range("J122").Value = 3289.5715
$ws.Range("L122").Value = 9868.7145
$ws.Range("N122").Value = -14768.7145
$ws.Range("H126").Value = 7435.143
$ws.Range("I126").Value = 7674.3335
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 23023.0005
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -20553.0005
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 3971.25
$ws.Range("I132").Value = 3971.25
$ws.Range("K132").Value = 11913.75
$ws.Range("M132").Value = -9383.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12981.429
$ws.Range("I7").Value = 13645
$ws.Range("K7").Value = 13645
$ws.Range("M7").Value = -13533
$ws.Range("H93").Value = 4000
$ws.Range("J93").Value = 4000
$ws.Range("L93").Value = 4000
$ws.Range("N93").Value = -6496
$ws.Range("H126").Value = 12981.429
$ws.Range("I126").Value = 13645
$ws.Range("K126").Value = 40935
$ws.Range("M126").Value = -38465
$ws.Range("H136").Value = 3114.111
$ws.Range("I136").Value = 2888.5386
$ws.Range("J136").Value = 3700.6
$ws.Range("K136").Value = 8665.6158
$ws.Range("L136").Value = 11101.8
$ws.Range("M136").Value = -6115.6158
$ws.Range("N136").Value = -16201.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3979.6667
$ws.Range("I62").Value = 3475.6
$ws.Range("K62").Value = 3475.6
$ws.Range("M62").Value = -2851.6
$ws.Range("H65").Value = 3979.6667
$ws.Range("I65").Value = 3475.6
$ws.Range("K65").Value = 17378
$ws.Range("M65").Value = -14258
$ws.Range("H126").Value = 33296.625
$ws.Range("I126").Value = 28849.8
$ws.Range("K126").Value = 86549.39999999999
$ws.Range("M126").Value = -84079.39999999999
$ws.Range("H132").Value = 1800.25
$ws.Range("I132").Value = 1800.25
$ws.Range("K132").Value = 5400.75
$ws.Range("M132").Value = -2870.75
